# Updates cryptos list (price/volume refresh) per the Dec 31 2023 GitHub
# Actions run. Column D (Price) / E (Volume(1h)) are plain text cells in
# the source sheet, so values that look numeric are entered with a
# leading apostrophe to keep them stored as text, then restored to the
# "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.892.99"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "2.290.33"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'314.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'105.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'39.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").Value = "'15.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "2.639.46"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "2.331.96"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "42.791.42"
$ws.Range("D19").Value = "'7.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'13.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +22.39%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'74.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'3.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").Value = "'10.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'7.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +23.20%  "
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "'22.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").Value = "'37.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'167.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("D35").Value = "'2.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("E36").Value = "  -4.27%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "'0.0352"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").Value = "'3.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("E41").Value = "  +4.76%  "
$ws.Range("D42").Value = "'70.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("D43").Value = "'0.233"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "'94.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "'12.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Value = "1.745.73"
$ws.Range("E47").Value = "  +9.90%  "
$ws.Range("D48").Value = "'113.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'80.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'5.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'8.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.42%  "
